# Rename the sheet from "Simple Template" to "Template"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Template"

# Adjust column widths (character-count units; Excel stores these internally
# with its own 1/256-character padding, which round-trips back to the
# target "x.83203125" style widths seen in the saved workbook)
$ws.Columns.Item(1).ColumnWidth = 50
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 15

# Update the filename value in A2
$ws.Range("A2").Value = "agent-261-17027502083-444.mp3"
